$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Withdraw History")

# Values that look numeric/date need a quote-prefix so Excel stores them
# as literal text (matching the existing column contents, e.g. "0900",
# "12/06/2021") rather than auto-converting to a number/date serial.
$ws.Range("A9").Value = "'200"
$ws.Range("B9").Value = "'0900"
$ws.Range("C9").Value = "'12/06/2021"
$ws.Range("D9").Value = "Lakeland, Florida"
$ws.Range("E9").Value = "N/A"

# Drop the quote-prefix formatting flag picked up above so the new cells
# don't end up with a distinct cell style from the rest of the sheet.
$ws.Range("A9:C9").ClearFormats()
